$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text values
$ws.Range("C1").Value = "ID Камеры"
$ws.Range("D1").Value = "Дата и время события `n гггг.мм.дд чч:мм:сс"
$ws.Range("E1").Value = "Распознанный ГРЗ"
$ws.Range("F1").Value = "Тип авто"

# Format header row: taller row, centered (no wrap), blue fill, thin black border
$headerRange = $ws.Range("A1:F1")
$ws.Rows(1).RowHeight = 40

$headerRange.HorizontalAlignment = -4108   # xlCenter
$headerRange.VerticalAlignment = -4108     # xlCenter
$headerRange.WrapText = $false

$headerRange.Interior.Color = 11824185     # OLE BGR value for RGB 39,6C,B4

$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2
$headerRange.Borders.Color = 0
